$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (raw OOXML stored width = ColumnWidth + 0.83) ---
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(31).ColumnWidth = 6.17

# --- Update data rows 2-5 with refreshed sensor readings (new simulation batch) ---
$row2 = New-Object 'object[,]' 1,34
$row2[0,0] = 45139.50694444445
$row2[0,1] = 19.217
$row2[0,2] = 12.901
$row2[0,3] = 4.042
$row2[0,4] = 40.812
$row2[0,5] = 32.818
$row2[0,6] = 15.123
$row2[0,7] = 47.986
$row2[0,8] = 23.269
$row2[0,9] = 9.710000000000001
$row2[0,10] = 14.67
$row2[0,11] = 16.076
$row2[0,12] = 16.742
$row2[0,13] = 4.827
$row2[0,14] = 15.038
$row2[0,15] = 20.994
$row2[0,16] = 12.85
$row2[0,17] = 3.46
$row2[0,18] = 2.249
$row2[0,19] = 221.547
$row2[0,20] = 41.81
$row2[0,21] = 13.881
$row2[0,22] = 27.553
$row2[0,23] = 14.055
$row2[0,24] = 3.03
$row2[0,25] = 24.312
$row2[0,26] = 12.261
$row2[0,27] = 11.125
$row2[0,28] = 13.047
$row2[0,29] = 16.565
$row2[0,30] = 3.456
$row2[0,31] = 42.557
$row2[0,32] = 7.647
$row2[0,33] = 17.354
$ws.Range("A2:AH2").Value = $row2

$row3 = New-Object 'object[,]' 1,34
$row3[0,0] = 45139.51388888889
$row3[0,1] = 21.139
$row3[0,2] = 15.167
$row3[0,3] = 1.968
$row3[0,4] = 45.737
$row3[0,5] = 37.222
$row3[0,6] = 16.635
$row3[0,7] = 63.468
$row3[0,8] = 25.596
$row3[0,9] = 11.164
$row3[0,10] = 16.545
$row3[0,11] = 18.312
$row3[0,12] = 19.19
$row3[0,13] = 5.314
$row3[0,14] = 16.542
$row3[0,15] = 23.396
$row3[0,16] = 14.113
$row3[0,17] = 1.595
$row3[0,18] = 1.204
$row3[0,19] = 244.48
$row3[0,20] = 46.269
$row3[0,21] = 15.269
$row3[0,22] = 30.859
$row3[0,23] = 16.077
$row3[0,24] = 2.746
$row3[0,25] = 31.003
$row3[0,26] = 13.487
$row3[0,27] = 12.108
$row3[0,28] = 14.21
$row3[0,29] = 19.056
$row3[0,30] = 1.265
$row3[0,31] = 57.69
$row3[0,32] = 8.516999999999999
$row3[0,33] = 19.09
$ws.Range("A3:AH3").Value = $row3

$row4 = New-Object 'object[,]' 1,34
$row4[0,0] = 45139.52083333334
$row4[0,1] = 20.658
$row4[0,2] = 15.056
$row4[0,3] = 1.463
$row4[0,4] = 44.824
$row4[0,5] = 36.609
$row4[0,6] = 16.258
$row4[0,7] = 63.625
$row4[0,8] = 25.014
$row4[0,9] = 11.023
$row4[0,10] = 16.31
$row4[0,11] = 17.985
$row4[0,12] = 18.897
$row4[0,13] = 5.193
$row4[0,14] = 16.166
$row4[0,15] = 22.944
$row4[0,16] = 13.716
$row4[0,17] = 1.087
$row4[0,18] = 0.9429999999999999
$row4[0,19] = 238.76
$row4[0,20] = 45.233
$row4[0,21] = 14.922
$row4[0,22] = 30.29
$row4[0,23] = 15.84
$row4[0,24] = 2.519
$row4[0,25] = 30.762
$row4[0,26] = 13.181
$row4[0,27] = 11.76
$row4[0,28] = 13.811
$row4[0,29] = 18.787
$row4[0,30] = 0.773
$row4[0,31] = 57.742
$row4[0,32] = 8.362
$row4[0,33] = 18.656
$ws.Range("A4:AH4").Value = $row4

$row5 = New-Object 'object[,]' 1,34
$row5[0,0] = 45139.52777777778
$row5[0,1] = 24.5
$row5[0,2] = 18.06
$row5[0,3] = 1.38
$row5[0,4] = 53.23
$row5[0,5] = 43.65
$row5[0,6] = 19.28
$row5[0,7] = 74.79000000000001
$row5[0,8] = 29.67
$row5[0,9] = 13.17
$row5[0,10] = 19.53
$row5[0,11] = 21.37
$row5[0,12] = 22.51
$row5[0,13] = 6.16
$row5[0,14] = 19.17
$row5[0,15] = 27.29
$row5[0,16] = 16.15
$row5[0,17] = 0.88
$row5[0,18] = 0.95
$row5[0,19] = 284.54
$row5[0,20] = 53.63
$row5[0,21] = 17.7
$row5[0,22] = 36.06
$row5[0,23] = 18.89
$row5[0,24] = 2.82
$row5[0,25] = 36.34
$row5[0,26] = 15.63
$row5[0,27] = 13.87
$row5[0,28] = 16.3
$row5[0,29] = 22.38
$row5[0,30] = 0.5600000000000001
$row5[0,31] = 67.86
$row5[0,32] = 9.970000000000001
$row5[0,33] = 22.13
$ws.Range("A5:AH5").Value = $row5

# --- Remove the now-obsolete last row (row 6) ---
$ws.Rows.Item(6).Delete()
